{"js": "// Add a \"Feedback\" section (heading + contact paragraph with a mailto\n// hyperlink) at the end of the document, wrapped in its own bookmark, and\n// normalize the \"Cabinet Office maintains...\" paragraph to a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 1. Merge the three runs of the \"Cabinet Office maintains...\" paragraph\n//    into a single run with the same combined text.\nlet cabinetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"The Cabinet Office maintains the following list of core descriptors\"\n    ) === 0\n  ) {\n    cabinetPara = paragraphs.items[i];\n    break;\n  }\n}\nif (cabinetPara) {\n  cabinetPara.clear();\n  await context.sync();\n  cabinetPara.insertText(\n    \"The Cabinet Office maintains the following list of core descriptors to ensure a consistent approach is adopted across all departments:\",\n    Word.InsertLocation.start\n  );\n  await context.sync();\n}\n\n// 2. Append the new \"Feedback\" heading and the contact paragraph after the\n//    last paragraph in the body (\"Descriptors are not codewords.\").\nconst lastParagraph = body.paragraphs.getLast();\n\nconst headingPara = lastParagraph.insertParagraph(\n  \"Feedback\",\n  Word.InsertLocation.after\n);\nheadingPara.style = \"Heading 2\";\nawait context.sync();\n\nconst contactPara = headingPara.insertParagraph(\n  \"\",\n  Word.InsertLocation.after\n);\ncontactPara.style = \"Block Text\";\nawait context.sync();\n\n// Build the paragraph content run-by-run so the sentence, the separating\n// space, the hyperlinked email address and the trailing period stay as\n// distinct runs (matching how the hyperlink-styled run is isolated).\ncontactPara.insertText(\n  \"If you have any questions or comments about this guidance, such as suggestions for improvements, please contact:\",\n  Word.InsertLocation.start\n);\nawait context.sync();\ncontactPara.insertText(\" \", Word.InsertLocation.end);\nawait context.sync();\ncontactPara.insertText(\n  \"itpolicycontent@digital.justice.gov.uk\",\n  Word.InsertLocation.end\n);\nawait context.sync();\ncontactPara.insertText(\".\", Word.InsertLocation.end);\nawait context.sync();\n\n// Turn the email address text into a real mailto: hyperlink.\nconst emailResults = contactPara.search(\"itpolicycontent@digital.justice.gov.uk\", {\n  matchCase: true,\n});\nemailResults.load(\"items\");\nawait context.sync();\nemailResults.items[0].hyperlink =\n  \"mailto:itpolicycontent@digital.justice.gov.uk\";\nawait context.sync();\n\n// 3. Wrap the new heading + contact paragraph in a bookmark named\n//    \"ariaid-title5\", matching the doc's \"ariaid-titleN\" heading-bookmark\n//    convention (the outer \"ariaid-title1\" bookmark already wraps the rest\n//    of the body, and its end simply moves past this new content).\nconst sectionRange = headingPara\n  .getRange(Word.RangeLocation.start)\n  .expandTo(contactPara.getRange(Word.RangeLocation.end));\nsectionRange.insertBookmark(\"ariaid-title5\");\nawait context.sync();\n", "ps1": "# Add a \"Feedback\" section (heading + contact paragraph with a mailto\n# hyperlink) at the end of the document, wrapped in its own bookmark, and\n# normalize the \"Cabinet Office maintains...\" paragraph to a single run.\n\n$d = $word.ActiveDocument\n\n# 1. Merge the three runs of the \"Cabinet Office maintains...\" paragraph\n#    into a single run with the same combined text.\n$cabinetText = \"The Cabinet Office maintains the following list of core descriptors to ensure a consistent approach is adopted across all departments:\"\n$cabinetPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.StartsWith(\"The Cabinet Office maintains\")) {\n        $cabinetPara = $d.Paragraphs($i)\n        break\n    }\n}\nif ($cabinetPara -ne $null) {\n    $cabinetRange = $cabinetPara.Range\n    $cabinetRange.MoveEnd(1, -1)  # exclude the paragraph mark\n    $cabinetRange.Delete()\n    $cabinetRange.InsertAfter($cabinetText)\n}\n\n# 2. Append the new \"Feedback\" heading and the contact paragraph after the\n#    last paragraph in the body (\"Descriptors are not codewords.\").\n$endRange = $d.Paragraphs.Last.Range\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.InsertAfter(\"Feedback\")\n$headingPara = $d.Paragraphs.Last\n$headingPara.Style = \"Heading 2\"\n\n$bodyRange = $headingPara.Range\n$bodyRange.Collapse(0)\n$bodyRange.InsertParagraphAfter()\n$bodyRange.Collapse(0)\n$bodyRange.InsertAfter(\"If you have any questions or comments about this guidance, such as suggestions for improvements, please contact: itpolicycontent@digital.justice.gov.uk.\")\n$contactPara = $d.Paragraphs.Last\n$contactPara.Style = \"Block Text\"\n\n# Turn the email address text into a real mailto: hyperlink.\n$findRange = $contactPara.Range.Duplicate\n$findRange.Find.Execute(\"itpolicycontent@digital.justice.gov.uk\")\n$d.Hyperlinks.Add($findRange, \"mailto:itpolicycontent@digital.justice.gov.uk\")\n\n# 3. Wrap the new heading + contact paragraph in a bookmark named\n#    \"ariaid-title5\", matching the doc's \"ariaid-titleN\" heading-bookmark\n#    convention (the outer \"ariaid-title1\" bookmark already wraps the rest\n#    of the body, and its end simply moves past this new content).\n$sectionRange = $d.Range($headingPara.Range.Start, $contactPara.Range.End)\n$d.Bookmarks.Add(\"ariaid-title5\", $sectionRange)\n"}
